$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Step A: Remove the duplicate bold "Play Diego Wild for Free..." title
# paragraph that sits right before the italic summary paragraph near the
# end of the document (it is being replaced by a proper "Meta
# description" paragraph right under the real H1 heading - see Step C).
# The same text also occurs once at the very top (the real H1 heading),
# so walk forward through every match and keep deleting the duplicate(s)
# that are NOT the first (heading) occurrence.
# -----------------------------------------------------------------------
$headingText = "Play Diego Wild for Free: See Gameplay, Payouts & Bonuses"
$matches = New-Object System.Collections.ArrayList
$searchRange = $d.Content
while ($searchRange.Find.Execute($headingText, `
        $true, $false, $false, $false, $false, $true, 0, $false, "", 0)) {
    $matches.Add(@($searchRange.Start, $searchRange.End)) | Out-Null
    $searchRange.Collapse(0)
}

# Delete every match EXCEPT the first one (the real H1 heading), walking
# backwards so earlier offsets stay valid. Expand each match to the full
# paragraph (including its trailing paragraph mark) before deleting so no
# empty paragraph is left behind.
for ($i = $matches.Count - 1; $i -ge 1; $i--) {
    $dupRange = $d.Range($matches[$i][0], $matches[$i][1])
    $dupRange.Expand(4) | Out-Null
    $dupRange.Delete()
}

# -----------------------------------------------------------------------
# Step B: Replace the text of the trailing italic paragraph with the new
# DALLE-prompt request text (formatting/italics is preserved because we
# use Find & Replace rather than rewriting the whole paragraph range).
# -----------------------------------------------------------------------
$oldBlurb = "Find out more about Diego Wild, the slot game set in the Amazon forest. Play for free or real money and explore ancient Aztec temples."
$newBlurb = "Could you please provide a detailed prompt for DALLE to create a feature image fitting the game " + [char]34 + "Diego Wild" + [char]34 + "? Prompt: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses - The Maya warrior should be surrounded by the game symbols, such as the bracelet, the vase, the statue, and the chest - The background of the image should represent the Amazon forest, with tall plants popping up on the sides of the game grid, arranged at the entrance to the ancient Aztec temple. - The colors of the image should be bright and eye-catching to grab the viewer's attention. Overall, the feature image should reflect the adventure/explorer theme of the game and convey a sense of excitement and fun. It should also showcase the main character, Diego, as a happy and playful warrior who is ready to embark on a thrilling adventure through the Amazon forest."

$d.Content.Find.Execute($oldBlurb, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newBlurb, 2) | Out-Null

# -----------------------------------------------------------------------
# Step C: Insert a new "Meta description" paragraph right after the H1
# title paragraph ("Play Diego Wild for Free: See Gameplay, Payouts &
# Bonuses"). The new paragraph has: an empty leading run, a bold run
# containing "Meta description", and a plain run with the rest of the
# meta-description sentence.
# -----------------------------------------------------------------------
$headingRange = $d.Content
$headingRange.Find.Execute("Play Diego Wild for Free: See Gameplay, Payouts & Bonuses", `
    $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$headingRange.InsertParagraphAfter() | Out-Null

$metaPara = $headingRange.Paragraphs.Item(1).Next()
$metaPara.Style = "Normal"

$boldLabel = "Meta description"
$restOfText = ": Find out more about Diego Wild, the slot game set in the Amazon forest. Play for free or real money and explore ancient Aztec temples."
$metaPara.Range.Text = $boldLabel + $restOfText

$metaStart = $metaPara.Range.Start
$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Font.Bold = 1
